# Shop.xlsx: "unify the conception of DataNode, DataTable, Entity."
#
# The author renamed the single worksheet from the stale "Property1" label
# to "DataNode" (matching the DataNode/DataTable/Entity naming used across
# the rest of the data-config workbooks), and left the cursor parked on
# D40 when they saved.
#
# (Everything else in the authoring diff - fileVersion/rupBuild bump,
# absPath, workbookView xWindow/yWindow/window size, xr:/xr2:/xr9: revision
# GUIDs, the MDW column-width/row-height re-rounding, the phonetic-guide
# font + phoneticPr, and the "Normal" -> "常规" cell-style locale label -
# are byproducts of the file having been re-saved by a different
# Office build/locale on Windows, not actions a script can take through
# the Excel object model, so they're intentionally left alone here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab.
$ws.Name = "DataNode"

# Leave the selection where the author left it before saving.
$ws.Range("D40").Select()
